$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# --- Overview sheet ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("B3").Value = $status
$ws1.Range("C3").Value = $status

# --- zh-cn sheet ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("B3").Value = $status
$ws2.Range("G3").Value = "2016-01-18 02:17:32"

# --- de-de sheet ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("B3").Value = $status
$ws3.Range("G3").Value = "2016-01-18 02:17:52"
